# Apply the diff to the "Truth Table" workbook:
#  - Set T20:T22 to 1 and T23:T25 to 0 (numeric values, replacing the "x" placeholder strings)
#  - Scroll/selection state: topLeftCell L1, active cell T26
#  - Window height change (cosmetic view size) on the workbook window

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T20").Value = 1
$ws.Range("T21").Value = 1
$ws.Range("T22").Value = 1
$ws.Range("T23").Value = 0
$ws.Range("T24").Value = 0
$ws.Range("T25").Value = 0

$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("T26").Select()

$excel.ActiveWindow.Height = 16620
